$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data table for Horeco (IGCC Netting Flows) - shift dates one day forward
# and refresh the import/export flow values.

$ws.Range("A2").Value = 46074
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 16.947
$ws.Range("A3").Value = 46074.01041666666
$ws.Range("B3").Value = 0.001
$ws.Range("C3").Value = 15.856
$ws.Range("A4").Value = 46074.02083333334
$ws.Range("B4").Value = 2.005
$ws.Range("C4").Value = 0.5580000000000001
$ws.Range("A5").Value = 46074.03125
$ws.Range("B5").Value = 5.714
$ws.Range("C5").Value = 0.021
$ws.Range("A6").Value = 46074.04166666666
$ws.Range("B6").Value = 3.562
$ws.Range("C6").Value = 0.425
$ws.Range("A7").Value = 46074.05208333334
$ws.Range("B7").Value = 5.297
$ws.Range("C7").Value = 0.856
$ws.Range("A8").Value = 46074.0625
$ws.Range("B8").Value = 6.691
$ws.Range("C8").Value = 0.533
$ws.Range("A9").Value = 46074.07291666666
$ws.Range("B9").Value = 3.756
$ws.Range("C9").Value = 0.005
$ws.Range("A10").Value = 46074.08333333334
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 4.734
$ws.Range("A11").Value = 46074.09375
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 2.018
$ws.Range("A12").Value = 46074.10416666666
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 4.442
$ws.Range("A13").Value = 46074.11458333334
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 3.24
$ws.Range("A14").Value = 46074.125
$ws.Range("B14").Value = 0.205
$ws.Range("C14").Value = 7.666
$ws.Range("A15").Value = 46074.13541666666
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 5.9
$ws.Range("A16").Value = 46074.14583333334
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 3.52
$ws.Range("A17").Value = 46074.15625
$ws.Range("B17").Value = 1.047
$ws.Range("C17").Value = 0.762
$ws.Range("A18").Value = 46074.16666666666
$ws.Range("B18").Value = 10.8
$ws.Range("C18").Value = 0
$ws.Range("A19").Value = 46074.17708333334
$ws.Range("B19").Value = 22.235
$ws.Range("C19").Value = 0
$ws.Range("A20").Value = 46074.1875
$ws.Range("B20").Value = 36.961
$ws.Range("C20").Value = 0
$ws.Range("A21").Value = 46074.19791666666
$ws.Range("B21").Value = 65.893
$ws.Range("C21").Value = 0
$ws.Range("A22").Value = 46074.20833333334
$ws.Range("B22").Value = 26.968
$ws.Range("C22").Value = 0
$ws.Range("A23").Value = 46074.21875
$ws.Range("B23").Value = 17.222
$ws.Range("C23").Value = 0
$ws.Range("A24").Value = 46074.22916666666
$ws.Range("B24").Value = 3.754
$ws.Range("C24").Value = 0.384
$ws.Range("A25").Value = 46074.23958333334
$ws.Range("B25").Value = 2.345
$ws.Range("C25").Value = 1.242
$ws.Range("A26").Value = 46074.25
$ws.Range("B26").Value = 1.629
$ws.Range("C26").Value = 0.854
$ws.Range("A27").Value = 46074.26041666666
$ws.Range("B27").Value = 29.084
$ws.Range("C27").Value = 0
$ws.Range("A28").Value = 46074.27083333334
$ws.Range("B28").Value = 29.924
$ws.Range("C28").Value = 0
$ws.Range("A29").Value = 46074.28125
$ws.Range("B29").Value = 66.84699999999999
$ws.Range("C29").Value = 0
$ws.Range("A30").Value = 46074.29166666666
$ws.Range("B30").Value = 62.704
$ws.Range("C30").Value = 0
$ws.Range("A31").Value = 46074.3125
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
